$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the rows for people no longer in the padron (bottom-to-top so
#    row numbers above stay valid while deleting):
#      row 35 -> UGO RUIZ            (Folio 3717)
#      row 31 -> PEPE DAVALOS        (Folio 3713)
#      row 21 -> JOSE LUIS COINDREAU (Folio 3703)
#      row 6  -> BAJA                (Folio 3688)
$ws.Rows.Item(35).EntireRow.Delete()
$ws.Rows.Item(31).EntireRow.Delete()
$ws.Rows.Item(21).EntireRow.Delete()
$ws.Rows.Item(6).EntireRow.Delete()

# 2. Rename "ALVARO SUAREZ X" -> "ALVARO SUAREZ" (still row 5 after the
#    above deletions, since all deleted rows were below it).
$ws.Range("B5").Value = "ALVARO SUAREZ"

# 2b. The new padron lists "SUSY 2" right after "SUSY", ahead of
#     "VICTOR FUENTES" (that pair is swapped vs. the old sheet's order).
#     After the deletions above, VICTOR FUENTES sits on row 32 and
#     SUSY 2 on row 33 - swap their names so SUSY 2 comes first.
$ws.Range("B32").Value = "SUSY 2"
$ws.Range("B33").Value = "VICTOR FUENTES"

# 3. Renumber the Folio column (A) starting at 3568 for the new padron.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$folio = 3568
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $folio
    $folio = $folio + 1
}

# 4. Re-style the data rows: plain font (size 11, black, Calibri), no
#    fill, no border - replacing the old banded/bordered look.
$dataRange = $ws.Range("A2:B" + $lastRow)
$dataRange.Style = "Normal"
$seed = $ws.Range("A2")
$seed.Font.Size = 11
$seed.Font.Color = 0
$seed.Copy()
$dataRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5. Move the selection like the saved workbook (selection moved to D47).
$ws.Range("D47").Select()

Write-Host "done"
